$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.935.58"
$ws.Range("D3").Value = "'1.624.25"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'213.66"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D9").Value = "'0.0615"
$ws.Range("E9").Value = "  -3.53%  "
$ws.Range("D10").Value = "'18.21"
$ws.Range("E10").Value = "  -6.76%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "'1.849.51"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.18"
$ws.Range("E13").Value = "  -1.95%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.619.74"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("E15").Value = "  -3.64%  "
$ws.Range("D16").Value = "'25.920.30"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").Value = "'61.18"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").Value = "'0.0₃0733"
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'191.57"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("E22").Value = "  -3.69%  "
$ws.Range("E23").Value = "  -2.29%  "
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'143.61"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("E28").Value = "  -2.43%  "
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  -3.98%  "
$ws.Range("D33").Value = "'3.10"
$ws.Range("E33").Value = "  -5.41%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").Value = "'1.117.95"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "'0.845"
$ws.Range("E37").Value = "  -6.45%  "
$ws.Range("E38").Value = "  -1.32%  "
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'97.93"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "'0.768"
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("D43").Value = "'1.760.54"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("E44").Value = "  -5.65%  "
$ws.Range("D45").Value = "'0.0₆0114"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("D47").Value = "'54.44"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  -3.72%  "
